# Update attendance/visitor counts (column F) across the 展览, 演出 and 全部类型
# sheets, per the upstream data refresh ("Update gh-pages to output generated
# at 456a3b4"). 本地生活 sheet is untouched by this refresh.

function Set-ColF($ws, $ref, $newVal) {
    $ws.Range($ref).Value = $newVal
}

$wb = $excel.ActiveWorkbook
$wsExpo = $wb.Worksheets.Item("展览")
$wsShow = $wb.Worksheets.Item("演出")
$wsAll  = $wb.Worksheets.Item("全部类型")

# 展览 sheet
Set-ColF $wsExpo "F6" 1396   # was 1393
Set-ColF $wsExpo "F8" 1104   # was 1103
Set-ColF $wsExpo "F9" 225   # was 221
Set-ColF $wsExpo "F10" 169   # was 165
Set-ColF $wsExpo "F11" 275   # was 274
Set-ColF $wsExpo "F12" 1732   # was 1728
Set-ColF $wsExpo "F13" 640   # was 639
Set-ColF $wsExpo "F14" 296   # was 293
Set-ColF $wsExpo "F15" 370   # was 369
Set-ColF $wsExpo "F16" 3880   # was 3859
Set-ColF $wsExpo "F20" 958   # was 957
Set-ColF $wsExpo "F21" 1243   # was 1242
Set-ColF $wsExpo "F24" 2905   # was 2904
Set-ColF $wsExpo "F25" 1730   # was 1727
Set-ColF $wsExpo "F30" 892   # was 891
Set-ColF $wsExpo "F31" 17   # was 16
Set-ColF $wsExpo "F32" 2107   # was 2103
Set-ColF $wsExpo "F33" 948   # was 946
Set-ColF $wsExpo "F34" 2165   # was 2153
Set-ColF $wsExpo "F36" 534   # was 530
Set-ColF $wsExpo "F37" 343   # was 325
Set-ColF $wsExpo "F41" 1008   # was 1002
Set-ColF $wsExpo "F42" 846   # was 845
Set-ColF $wsExpo "F43" 1113   # was 1111
Set-ColF $wsExpo "F44" 244   # was 241
Set-ColF $wsExpo "F45" 475   # was 473
Set-ColF $wsExpo "F46" 325   # was 323
Set-ColF $wsExpo "F48" 3418   # was 3416

# 演出 sheet
Set-ColF $wsShow "F11" 850   # was 849
Set-ColF $wsShow "F12" 28   # was 27
Set-ColF $wsShow "F16" 12   # was 11

# 全部类型 sheet
Set-ColF $wsAll "F5" 1396   # was 1393
Set-ColF $wsAll "F7" 1104   # was 1103
Set-ColF $wsAll "F8" 225   # was 221
Set-ColF $wsAll "F9" 169   # was 165
Set-ColF $wsAll "F11" 1732   # was 1728
Set-ColF $wsAll "F12" 640   # was 639
Set-ColF $wsAll "F13" 296   # was 293
Set-ColF $wsAll "F14" 370   # was 369
Set-ColF $wsAll "F15" 3880   # was 3859
Set-ColF $wsAll "F21" 1243   # was 1242
Set-ColF $wsAll "F22" 2905   # was 2904
Set-ColF $wsAll "F24" 1730   # was 1727
Set-ColF $wsAll "F29" 850   # was 849
Set-ColF $wsAll "F30" 28   # was 27
Set-ColF $wsAll "F32" 892   # was 891
Set-ColF $wsAll "F33" 2107   # was 2103
Set-ColF $wsAll "F34" 12   # was 11
Set-ColF $wsAll "F35" 948   # was 946
Set-ColF $wsAll "F36" 2166   # was 2153
Set-ColF $wsAll "F37" 534   # was 530
Set-ColF $wsAll "F38" 343   # was 325
Set-ColF $wsAll "F40" 1008   # was 1002
Set-ColF $wsAll "F41" 846   # was 845
Set-ColF $wsAll "F42" 1113   # was 1111
Set-ColF $wsAll "F43" 244   # was 242
Set-ColF $wsAll "F44" 475   # was 473
Set-ColF $wsAll "F45" 325   # was 323
Set-ColF $wsAll "F48" 3418   # was 3416

$wb.Application.CalculateFullRebuild()
